$d = $word.ActiveDocument

# 1) Rename the function everywhere it appears verbatim (bold heading run,
#    the "Run the tool <name>" mention further down, and as a substring of
#    the SQL example -- the SQL example gets extra handling below).
$d.Content.Find.Execute("ni_connect_hanging_edges_to_nodes", $false, $false, $false, $false, $false, $true, 1, $false, "ni_data_proc_connect_hanging_edges_to_nodes_in_search", 2)

# 2) In the SQL example a space was inserted between the (renamed) function
#    name and the opening paren. Locate the exact text and overwrite the
#    Range.Text directly (rather than via Find's replacement argument) so
#    the literal single-quote characters are not mangled into curly quotes
#    by the smart-quotes simulation.
$rng = $d.Content
$rng.Find.Execute("ni_data_proc_connect_hanging_edges_to_nodes_in_search('data_national_grid_gas_pipeline_feeder'", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "ni_data_proc_connect_hanging_edges_to_nodes_in_search ('data_national_grid_gas_pipeline_feeder'"

# 3) Parameter 1 description: Edge_table_prefix -> Edge_table_name
$rng = $d.Content
$rng.Find.Execute("Edge_table_prefix", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "Edge_table_name"

# 4) Parameter 4 description: Node_table_prefix -> Node_table_name
$rng = $d.Content
$rng.Find.Execute("Node_table_prefix", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "Node_table_name"

# 5) Final paragraph: "...newly derived geometry (bl)," -> "...(bl)."
$rng = $d.Content
$rng.Find.Execute("bl),", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$commaRng = $d.Range($rng.End - 1, $rng.End)
$commaRng.Text = "."

# 6) Append a new paragraph explaining the _join/_unique output tables.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRng = $last.Range
$endRng.Collapse(0)
$endRng.Text = "`rThis table would be output with _join appended to the input output table name (parameter 7). Secondly a table with _unique appended to the input output table name (parameter 7) is also written to the schema that contains the original geometry replaced with the newly derived geometry."
